$d = $word.ActiveDocument

function Replace-ExactText($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-ExactText "Kichwa cha Video" "Video Title"
Replace-ExactText "Mada" "Topic"
Replace-ExactText "Jiometri" "Geometry"
Replace-ExactText "Malengo" "Aim(s)"
Replace-ExactText "Urefu" "Length"
Replace-ExactText "Mahali pa Kambi" "Camp Location"
Replace-ExactText "Wawezeshaji" "Facilitators"
Replace-ExactText "N. ya wanafunzi" "N. of students"
Replace-ExactText "Tarehe" "Date"
Replace-ExactText "Rasilimali" "Resources"
Replace-ExactText "inahitajika" "needed"
Replace-ExactText "Maandalizi" "Preparations"
Replace-ExactText "Muda wa video" "Video time"
Replace-ExactText "Mwezeshaji anafanya nini" "What facilitator does"
Replace-ExactText "Wanachofanya wanafunzi" "What learners do"
Replace-ExactText "Utangulizi Mkuu wa Video ya VMC" "General VMC Video Introduction"

# Update the document default language from Swahili (Kenya) to Swahili (Tanzania)
$d.Styles("Normal").LanguageID = "sw-TZ"
